$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 5-7 (old FAPs/MuSCs -> MuSCs and MuSCs/ECs, MuSCs/MuSCs rows), shrinking the table to 3 data rows
$ws.Range("A5:T7").EntireRow.Delete()

# Update remaining data rows (2-4) with the refreshed TPM-derived values
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Tnfsf13b"
$ws.Range("C2").Value = "Tnfrsf13b"
$ws.Range("D2").Value = "MuSCs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.04408699999999999
$ws.Range("H2").Value = 0.132261
$ws.Range("I2").Value = 0.007006504090795892
$ws.Range("J2").Value = 0.007006504090795892
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.04819200000000001
$ws.Range("N2").Value = 0.144576
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.002124640704
$ws.Range("R2").Value = 0.019121766336
$ws.Range("S2").Value = 0.007006504090795892
$ws.Range("T2").Value = 0.007006504090795892
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Tnfsf13b"
$ws.Range("C3").Value = "Tnfrsf13b"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 6.019967333333334
$ws.Range("H3").Value = 18.059902
$ws.Range("I3").Value = 0.9567202519440571
$ws.Range("J3").Value = 0.9567202519440571
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.04819200000000001
$ws.Range("N3").Value = 0.144576
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.290114265728
$ws.Range("R3").Value = 2.611028391552
$ws.Range("S3").Value = 0.9567202519440571
$ws.Range("T3").Value = 0.9567202519440571
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Tnfsf13b"
$ws.Range("C4").Value = "Tnfrsf13b"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.228242
$ws.Range("H4").Value = 0.6847259999999999
$ws.Range("I4").Value = 0.03627324396514701
$ws.Range("J4").Value = 0.03627324396514701
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.04819200000000001
$ws.Range("N4").Value = 0.144576
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.010999438464
$ws.Range("R4").Value = 0.09899494617599999
$ws.Range("S4").Value = 0.03627324396514701
$ws.Range("T4").Value = 0.03627324396514701
